# FQ 9.1-001-01 Indicadores da Qualidade - roll the dashboard from 2024 to 2025
# and update the Jan/Fev actuals for "PAQ - Prazo e Devolução".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text: bump the year shown above the chart ---
$ws.Range("A1").Value = "Índice de Entrega Do Provedor Externo no Prazo - 2025"

# --- Main monthly table (rows 8-19): Jan/Fev actuals updated, Mar-Nov cleared ---
$ws.Range("B8").Value = 150
$ws.Range("C8").Value = 149
$ws.Range("H8").Value = 1

$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 100

$ws.Range("B10:C18").ClearContents()

# --- "TABELA GERAL DO PROVEDOR EXTERNO" helper table (rows 24-35), column B ---
$ws.Range("B24").Value = 150
$ws.Range("B25").Value = 100
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0
$ws.Range("B32").Value = 0
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = 0

# --- Devoluções helper table (rows 40-51), column B: Out (row 49) is a hardcoded value ---
$ws.Range("B49").Value = 0

# --- restore the on-screen selection to where the user left off ---
$ws.Range("B42").Select()
